# Auto-generated script applying scheduled-runner market data refresh
# to the Sheets workbook (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For every affected row, columns H/I/J/K/L/M/N (current market prices,
# leve sale prices, and computed profits) are updated to the latest
# refreshed values. A few cells whose profit became undefined/zero
# (denominator columns dropped to 0) are cleared entirely, matching
# source behavior of omitting such cells.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3448566.8
$ws.Range("J6").Value = 366.16666
$ws.Range("L6").Value = 1098.49998
$ws.Range("N6").Value = -1322.49998

$ws.Range("H76").Value = 14165.667
$ws.Range("I76").Value = 13399.4
$ws.Range("K76").Value = 13399.4
$ws.Range("M76").Value = -13084.4

$ws.Range("H79").Value = 14165.667
$ws.Range("I79").Value = 13399.4
$ws.Range("K79").Value = 13399.4
$ws.Range("M79").Value = -12307.4

$ws.Range("H96").Value = 683.5
$ws.Range("I96").Value = 119.5
$ws.Range("K96").Value = 358.5
$ws.Range("M96").Value = 1014.5

$ws.Range("H98").Value = 2100.658
$ws.Range("I98").Value = 2103.4324
$ws.Range("K98").Value = 2103.4324
$ws.Range("M98").Value = -605.4324000000001

$ws.Range("H103").Value = 1345.5
$ws.Range("I103").Value = 1239.5
$ws.Range("J103").Value = 1372
$ws.Range("K103").Value = 3718.5
$ws.Range("L103").Value = 4116
$ws.Range("M103").Value = -3132.5
$ws.Range("N103").Value = -5288

$ws.Range("H111").Value = 10418465
$ws.Range("I111").Value = 17857700
$ws.Range("K111").Value = 53573100
$ws.Range("M111").Value = -53570033

$ws.Range("H113").Value = 89519130
$ws.Range("J113").Value = 115401670
$ws.Range("L113").Value = 115401670
$ws.Range("N113").Value = -115408178

$ws.Range("H122").Value = 2100.658
$ws.Range("I122").Value = 2103.4324
$ws.Range("K122").Value = 6310.297200000001
$ws.Range("M122").Value = -3860.297200000001

$ws.Range("H132").Value = 748.46295
$ws.Range("I132").Value = 748.46295
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2245.38885
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 284.6111500000002
$ws.Range("N132").ClearContents()

$ws.Range("H137").Value = 6755.24
$ws.Range("I137").Value = 4249.3
$ws.Range("K137").Value = 12747.9
$ws.Range("M137").Value = -10197.9


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3285610.2
$ws.Range("I32").Value = 3395292
$ws.Range("K32").Value = 3395292
$ws.Range("M32").Value = -3395005

$ws.Range("H45").Value = 5108.1113
$ws.Range("I45").Value = 1697.4286
$ws.Range("K45").Value = 1697.4286
$ws.Range("M45").Value = -1320.4286

$ws.Range("H61").Value = 38467624
$ws.Range("I61").Value = 2175.875
$ws.Range("K61").Value = 2175.875
$ws.Range("M61").Value = -1963.875

$ws.Range("H74").Value = 24562.555
$ws.Range("I74").Value = 30171.629
$ws.Range("J74").Value = 4930.8
$ws.Range("K74").Value = 30171.629
$ws.Range("L74").Value = 4930.8
$ws.Range("M74").Value = -29297.629
$ws.Range("N74").Value = -6678.8

$ws.Range("H77").Value = 24562.555
$ws.Range("I77").Value = 30171.629
$ws.Range("J77").Value = 4930.8
$ws.Range("K77").Value = 150858.145
$ws.Range("L77").Value = 24654
$ws.Range("M77").Value = -146490.145
$ws.Range("N77").Value = -33390

$ws.Range("H110").Value = 55556804
$ws.Range("I110").Value = 1472.25
$ws.Range("J110").Value = 166667470
$ws.Range("K110").Value = 1472.25
$ws.Range("L110").Value = 166667470
$ws.Range("M110").Value = 572.75
$ws.Range("N110").Value = -166671560

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H132").Value = 4951.75
$ws.Range("I132").Value = 1190.625
$ws.Range("K132").Value = 3571.875
$ws.Range("M132").Value = -1041.875

$ws.Range("H136").Value = 38467624
$ws.Range("I136").Value = 2175.875
$ws.Range("K136").Value = 6527.625
$ws.Range("M136").Value = -3977.625


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 112501200
$ws.Range("I107").Value = 125001110
$ws.Range("K107").Value = 125001110
$ws.Range("M107").Value = -124999190

$ws.Range("I134").Value = 12501503
$ws.Range("J134").Value = 10416.471
$ws.Range("K134").Value = 37504509
$ws.Range("L134").Value = 31249.413
$ws.Range("M134").Value = -37501974
$ws.Range("N134").Value = -36319.413


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 44323.2
$ws.Range("J51").Value = 44729
$ws.Range("L51").Value = 44729
$ws.Range("N51").Value = -46201

$ws.Range("H58").Value = 6862.2163
$ws.Range("I58").Value = 2797.2307
$ws.Range("J58").Value = 9064.083000000001
$ws.Range("K58").Value = 2797.2307
$ws.Range("L58").Value = 9064.083000000001
$ws.Range("M58").Value = -2594.2307
$ws.Range("N58").Value = -9470.083000000001

$ws.Range("H61").Value = 44323.2
$ws.Range("J61").Value = 44729
$ws.Range("L61").Value = 44729
$ws.Range("N61").Value = -45425

$ws.Range("H135").Value = 64000
$ws.Range("I135").Value = 64000
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 64000
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -58930
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 6862.2163
$ws.Range("I136").Value = 2797.2307
$ws.Range("J136").Value = 9064.083000000001
$ws.Range("K136").Value = 8391.6921
$ws.Range("L136").Value = 27192.249
$ws.Range("M136").Value = -5841.6921
$ws.Range("N136").Value = -32292.249


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 67903450
$ws.Range("I109").Value = 62502490
$ws.Range("K109").Value = 187507470
$ws.Range("M109").Value = -187506430

$ws.Range("H113").Value = 4618.04
$ws.Range("J113").Value = 6129.2354
$ws.Range("L113").Value = 18387.7062
$ws.Range("N113").Value = -22727.7062

$ws.Range("H132").Value = 13029.591
$ws.Range("I132").Value = 3905.3
$ws.Range("J132").Value = 20633.166
$ws.Range("K132").Value = 35147.7
$ws.Range("L132").Value = 185698.494
$ws.Range("M132").Value = -32617.7
$ws.Range("N132").Value = -190758.494

$ws.Range("H134").Value = 49948.59
$ws.Range("I134").Value = 54118.45
$ws.Range("J134").Value = 8250
$ws.Range("K134").Value = 162355.35
$ws.Range("L134").Value = 24750
$ws.Range("M134").Value = -157285.35
$ws.Range("N134").Value = -34890

$ws.Range("H136").Value = 20835348
$ws.Range("I136").Value = 20835348
$ws.Range("K136").Value = 62506044
$ws.Range("M136").Value = -62500944

$ws.Range("H137").Value = 90305.56
$ws.Range("I137").Value = 74259.5
$ws.Range("J137").Value = 115266.11
$ws.Range("K137").Value = 222778.5
$ws.Range("L137").Value = 345798.33
$ws.Range("M137").Value = -217678.5
$ws.Range("N137").Value = -355998.33


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2920.5833
$ws.Range("I80").Value = 2791.5
$ws.Range("K80").Value = 2791.5
$ws.Range("M80").Value = -1793.5

$ws.Range("H83").Value = 2920.5833
$ws.Range("I83").Value = 2791.5
$ws.Range("K83").Value = 13957.5
$ws.Range("M83").Value = -8965.5

$ws.Range("H102").Value = 2169.45
$ws.Range("I102").Value = 2234
$ws.Range("J102").Value = 1803.6666
$ws.Range("K102").Value = 2234
$ws.Range("L102").Value = 1803.6666
$ws.Range("M102").Value = -612
$ws.Range("N102").Value = -5047.6666

$ws.Range("H132").Value = 4056.8408
$ws.Range("I132").Value = 2219.75
$ws.Range("J132").Value = 8955.75
$ws.Range("K132").Value = 6659.25
$ws.Range("L132").Value = 26867.25
$ws.Range("M132").Value = -4129.25
$ws.Range("N132").Value = -31927.25

$ws.Range("H136").Value = 28772.193
$ws.Range("J136").Value = 29343.77
$ws.Range("L136").Value = 88031.31
$ws.Range("N136").Value = -93131.31


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4245.4585
$ws.Range("I40").Value = 3139.2
$ws.Range("K40").Value = 3139.2
$ws.Range("M40").Value = -3003.2

$ws.Range("H136").Value = 10758.259
$ws.Range("I136").Value = 2903.92
$ws.Range("J136").Value = 17529.242
$ws.Range("K136").Value = 8711.76
$ws.Range("L136").Value = 52587.726
$ws.Range("M136").Value = -6161.76
$ws.Range("N136").Value = -57687.726


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1046.4
$ws.Range("I96").Value = 999.5
$ws.Range("K96").Value = 999.5
$ws.Range("M96").Value = 373.5

$ws.Range("H100").Value = 531.4643
$ws.Range("I100").Value = 350.44446
$ws.Range("J100").Value = 857.3
$ws.Range("K100").Value = 700.88892
$ws.Range("L100").Value = 1714.6
$ws.Range("M100").Value = -159.88892
$ws.Range("N100").Value = -2796.6

$ws.Range("H122").Value = 225882.22
$ws.Range("J122").Value = 7281
$ws.Range("L122").Value = 21843
$ws.Range("N122").Value = -26743

$ws.Range("H132").Value = 7082.517
$ws.Range("J132").Value = 5782.5454
$ws.Range("L132").Value = 17347.6362
$ws.Range("N132").Value = -22407.6362

